$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Testcase sample IDs shift down: the previous A2 value moves to A3, and a
# new sample ID is recorded in A2.
$ws.Range("A3").Value = "A0732302"
$ws.Range("A2").Value = "A1286705789"
